$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# Title Placeholder 1: "Click to edit Master title style"
$m.Shapes.Item(1).TextFrame.TextRange.Text = "{g0}ickclay otay edithay astermay itletay estylay{/g1}"

# Text Placeholder 2: 5 outline-level paragraphs
$bodyTr = $m.Shapes.Item(2).TextFrame.TextRange
$bodyTr.Paragraphs(1,1).Text = "{g0}ickclay otay edithay astermay exttay esstylay{/g1}"
$bodyTr.Paragraphs(2,1).Text = "{g0}econdsay evellay{/g1}"
$bodyTr.Paragraphs(3,1).Text = "{g0}irdthay evellay{/g1}"
$bodyTr.Paragraphs(4,1).Text = "{g0}ourthfay evellay{/g1}"
$bodyTr.Paragraphs(5,1).Text = "{g0}ifthfay evellay{/g1}"

# Date Placeholder 3: datetimeFigureOut field text "3/1/2007"
$m.Shapes.Item(3).TextFrame.TextRange.Text = "{g0}3/1/2007{/g1}"

# Slide Number Placeholder 5: slidenum field text "<#>"
$m.Shapes.Item(5).TextFrame.TextRange.Text = "{g0}‹#›{/g1}"
